$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 44817
$ws.Range("J4").Value = 85
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15529
$ws.Range("P4").Value = 1035

$ws.Range("D5").Value = 44754
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("P5").Value = 1000

$ws.Range("D6").Value = 44313
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("P6").Value = 933

$ws.Range("D7").Value = 44742
$ws.Range("J7").Value = 85
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15529
$ws.Range("P7").Value = 1035

$ws.Range("D8").Value = 44397
$ws.Range("J8").Value = 73
$ws.Range("K8").Value = 21000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21521
$ws.Range("P8").Value = 1435

$ws.Range("D9").Value = 44333
$ws.Range("J9").Value = 35
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 15000
$ws.Range("P9").Value = 1000

$ws.Range("D10").Value = 44729
$ws.Range("J10").Value = 85
$ws.Range("K10").Value = 16000
$ws.Range("L10").Value = 17000
$ws.Range("M10").Value = 16529
$ws.Range("P10").Value = 1102

$ws.Range("D11").Value = 44762
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("P11").Value = 967

$ws.Range("D14").Value = 44329
$ws.Range("J14").Value = 35
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 15000
$ws.Range("P14").Value = 1000

$ws.Range("D15").Value = 44736
$ws.Range("J15").Value = 82
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 16488
$ws.Range("P15").Value = 1099

$ws.Range("D16").Value = 44334
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = 14000
$ws.Range("P16").Value = 933

$ws.Range("D17").Value = 44753
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 1033

$ws.Range("D18").Value = 44811
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 16000
$ws.Range("M18").Value = 16000
$ws.Range("P18").Value = 1067

$ws.Range("D19").Value = 44370
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 18000
$ws.Range("P19").Value = 1200

$ws.Range("D20").Value = 44750
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15471
$ws.Range("P20").Value = 1031

$ws.Range("D21").Value = 44769
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14471
$ws.Range("P21").Value = 965

$ws.Range("D22").Value = 44438
$ws.Range("J22").Value = 75
$ws.Range("K22").Value = 19000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 19467
$ws.Range("P22").Value = 1298

$ws.Range("D23").Value = 44725
$ws.Range("J23").Value = 85
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14471
$ws.Range("P23").Value = 965

$ws.Range("D24").Value = 44804
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("P24").Value = 1000

$ws.Range("D25").Value = 44720
$ws.Range("J25").Value = 85
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 15529
$ws.Range("P25").Value = 1035

$ws.Range("D26").Value = 44791
$ws.Range("J26").Value = 40
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 15000
$ws.Range("P26").Value = 1000

$ws.Range("D27").Value = 44748
$ws.Range("J27").Value = 73
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 15521
$ws.Range("P27").Value = 1035

$ws.Range("D28").Value = 44783
$ws.Range("J28").Value = 50
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 15000
$ws.Range("P28").Value = 1000

$ws.Range("D29").Value = 44722
$ws.Range("J29").Value = 95
$ws.Range("L29").Value = 15500
$ws.Range("M29").Value = 15263
$ws.Range("P29").Value = 1018

$ws.Range("D30").Value = 44312
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 13000
$ws.Range("L30").Value = 14000
$ws.Range("M30").Value = 13562
$ws.Range("P30").Value = 904

$ws.Range("D31").Value = 44756
$ws.Range("J31").Value = 50
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = 15000
$ws.Range("P31").Value = 1000

$ws.Range("D32").Value = 44797
$ws.Range("J32").Value = 40

$ws.Range("D33").Value = 44799
$ws.Range("J33").Value = 55

$ws.Range("D34").Value = 44792
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = 15000
$ws.Range("P34").Value = 1000

$ws.Range("D35").Value = 44803
$ws.Range("J35").Value = 85
$ws.Range("L35").Value = 15500
$ws.Range("M35").Value = 15265
$ws.Range("P35").Value = 1018

$ws.Range("D36").Value = 44755
$ws.Range("J36").Value = 100
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 15550
$ws.Range("P36").Value = 1037

$ws.Range("D37").Value = 44714
$ws.Range("J37").Value = 100
$ws.Range("M37").Value = 15250
$ws.Range("P37").Value = 1017

$ws.Range("D38").Value = 44330
$ws.Range("J38").Value = 30
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 15000
$ws.Range("P38").Value = 1000

$ws.Range("D39").Value = 44810
$ws.Range("J39").Value = 85
$ws.Range("K39").Value = 16000
$ws.Range("L39").Value = 16500
$ws.Range("M39").Value = 16235
$ws.Range("P39").Value = 1082

$ws.Range("D40").Value = 44789
$ws.Range("J40").Value = 40

$ws.Range("D41").Value = 44746
$ws.Range("J41").Value = 103
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = 15563
$ws.Range("P41").Value = 1038

$ws.Range("D42").Value = 44308
$ws.Range("K42").Value = 16000
$ws.Range("L42").Value = 16000
$ws.Range("M42").Value = 16000
$ws.Range("P42").Value = 1067

$ws.Range("D43").Value = 44340
$ws.Range("J43").Value = 47
$ws.Range("K43").Value = 14000
$ws.Range("L43").Value = 14000
$ws.Range("M43").Value = 14000
$ws.Range("P43").Value = 933

$ws.Range("D44").Value = 44323
$ws.Range("K44").Value = 15000
$ws.Range("L44").Value = 15000
$ws.Range("M44").Value = 15000
$ws.Range("P44").Value = 1000

$ws.Range("D45").Value = 44790
$ws.Range("J45").Value = 40
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = 15000
$ws.Range("P45").Value = 1000

$ws.Range("D46").Value = 44785
$ws.Range("J46").Value = 85
$ws.Range("K46").Value = 14000
$ws.Range("M46").Value = 14471
$ws.Range("P46").Value = 965

$ws.Range("D47").Value = 44721
$ws.Range("J47").Value = 130
$ws.Range("K47").Value = 14000
$ws.Range("M47").Value = 14500
$ws.Range("P47").Value = 967

$ws.Range("D48").Value = 44377
$ws.Range("J48").Value = 80
$ws.Range("K48").Value = 18000
$ws.Range("L48").Value = 19000
$ws.Range("M48").Value = 18500
$ws.Range("P48").Value = 1233

$ws.Range("D49").Value = 44784
$ws.Range("J49").Value = 105
$ws.Range("M49").Value = 14476
$ws.Range("P49").Value = 965

$ws.Range("D50").Value = 44818
$ws.Range("J50").Value = 58
$ws.Range("K50").Value = 16000
$ws.Range("L50").Value = 16000
$ws.Range("M50").Value = 16000
$ws.Range("P50").Value = 1067

$ws.Range("D51").Value = 44719
$ws.Range("J51").Value = 60
$ws.Range("K51").Value = 15000
$ws.Range("M51").Value = 15000
$ws.Range("P51").Value = 1000

$ws.Range("D52").Value = 44757
$ws.Range("J52").Value = 40

$ws.Range("D53").Value = 44767
$ws.Range("J53").Value = 45

$ws.Range("D54").Value = 44776
$ws.Range("J54").Value = 105
$ws.Range("L54").Value = 15500
$ws.Range("M54").Value = 15238
$ws.Range("P54").Value = 1016

$ws.Range("D55").Value = 44326
$ws.Range("J55").Value = 45
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = 15000
$ws.Range("P55").Value = 1000

$ws.Range("D56").Value = 44314

$ws.Range("D57").Value = 44341
$ws.Range("J57").Value = 40

$ws.Range("D58").Value = 44747
$ws.Range("K58").Value = 16000
$ws.Range("L58").Value = 16000
$ws.Range("M58").Value = 16000
$ws.Range("P58").Value = 1067

$ws.Range("D59").Value = 44806
$ws.Range("J59").Value = 45

$ws.Range("D60").Value = 44321
$ws.Range("J60").Value = 38
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = 15000
$ws.Range("P60").Value = 1000

$ws.Range("D61").Value = 44715
$ws.Range("J61").Value = 85
$ws.Range("L61").Value = 15500
$ws.Range("M61").Value = 15235
$ws.Range("P61").Value = 1016

$ws.Range("D62").Value = 44795
$ws.Range("J62").Value = 56
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = 15000
$ws.Range("P62").Value = 1000

$ws.Range("D63").Value = 44775
$ws.Range("J63").Value = 93
$ws.Range("K63").Value = 14000
$ws.Range("M63").Value = 14516
$ws.Range("P63").Value = 968

$ws.Range("D64").Value = 44448
$ws.Range("J64").Value = 85
$ws.Range("K64").Value = 21000
$ws.Range("L64").Value = 22000
$ws.Range("M64").Value = 21529
$ws.Range("P64").Value = 1435

$ws.Range("D65").Value = 44809
$ws.Range("J65").Value = 105
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = 15476
$ws.Range("P65").Value = 1032

$ws.Range("D66").Value = 44315
$ws.Range("J66").Value = 65
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = 14538
$ws.Range("P66").Value = 969

$ws.Range("D67").Value = 44327
$ws.Range("J67").Value = 35
$ws.Range("K67").Value = 15000
$ws.Range("M67").Value = 15000
$ws.Range("P67").Value = 1000

$ws.Range("D68").Value = 44316
$ws.Range("J68").Value = 45
$ws.Range("K68").Value = 14000
$ws.Range("M68").Value = 14444
$ws.Range("P68").Value = 963

$ws.Range("D69").Value = 44320
$ws.Range("J69").Value = 40
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = 15000
$ws.Range("P69").Value = 1000

$ws.Range("D70").Value = 44764
$ws.Range("J70").Value = 45

$ws.Range("D71").Value = 44749
$ws.Range("J71").Value = 100
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 17000
$ws.Range("M71").Value = 16450
$ws.Range("P71").Value = 1097

$ws.Range("D72").Value = 44727
$ws.Range("J72").Value = 60
$ws.Range("K72").Value = 15000
$ws.Range("L72").Value = 15000
$ws.Range("M72").Value = 15000
$ws.Range("P72").Value = 1000

$ws.Range("D73").Value = 44761
$ws.Range("J73").Value = 50

$ws.Range("D74").Value = 44798
$ws.Range("J74").Value = 100
$ws.Range("K74").Value = 14000
$ws.Range("M74").Value = 14450
$ws.Range("P74").Value = 963

$ws.Range("D75").Value = 44763
$ws.Range("J75").Value = 80
$ws.Range("M75").Value = 14500
$ws.Range("P75").Value = 967

$ws.Range("D76").Value = 44816
$ws.Range("J76").Value = 60
$ws.Range("K76").Value = 16000
$ws.Range("L76").Value = 16000
$ws.Range("M76").Value = 16000
$ws.Range("P76").Value = 1067

# Remove the last data row (was row 87), matching dimension change A1:R87 -> A1:R86
$ws.Rows(87).Delete()
